# Apply updated crypto price/volume figures to columns D (Price) and E (Volume(1h))
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.903.05"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").Value = "1.830.18"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("D4").Value = "'1.009"
$ws.Range("E4").Value = "  +0.78%  "
$ws.Range("D5").Value = "'311.13"
$ws.Range("E5").Value = "  -1.01%  "
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("D7").Value = "'0.4569"
$ws.Range("E7").Value = "  -0.97%  "
$ws.Range("D8").Value = "'0.3689"
$ws.Range("E8").Value = "  -0.42%  "
$ws.Range("D9").Value = "'0.07179"
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("D10").Value = "'0.8768"
$ws.Range("D11").Value = "'0.07844"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("D13").Value = "1.827.24"
$ws.Range("E13").Value = "  -2.11%  "
$ws.Range("E14").Value = "  -0.90%  "
$ws.Range("D15").Value = "'6.391"
$ws.Range("E15").Value = "  -2.56%  "
$ws.Range("D16").Value = "'87.14"
$ws.Range("E16").Value = "  -5.19%  "
$ws.Range("D17").Value = "'1.010"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").Value = "'0.000008722"
$ws.Range("E18").Value = "  -1.37%  "
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").Value = "26.932.42"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("E21").Value = "  -2.21%  "
$ws.Range("D22").Value = "'5.006"
$ws.Range("E22").Value = "  -2.23%  "
$ws.Range("D23").Value = "2.046.09"
$ws.Range("E23").Value = "  -3.73%  "
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("D25").Value = "'1.983"
$ws.Range("E25").Value = "  +5.06%  "
$ws.Range("D26").Value = "'151.26"
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("D27").Value = "'18.21"
$ws.Range("E27").Value = "  -0.92%  "
$ws.Range("D28").Value = "'1.969"
$ws.Range("E28").Value = "  -5.45%  "
$ws.Range("D29").Value = "'113.92"
$ws.Range("E29").Value = "  -1.80%  "
$ws.Range("D30").Value = "'4.930"
$ws.Range("E30").Value = "  -3.90%  "
$ws.Range("D31").Value = "'0.08799"
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("E32").Value = "  +1.11%  "
$ws.Range("D33").Value = "'0.7547"
$ws.Range("E33").Value = "  -0.86%  "
$ws.Range("D34").Value = "'4.483"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").Value = "'1.132"
$ws.Range("E35").Value = "  -3.40%  "
$ws.Range("D36").Value = "'2.571"
$ws.Range("E36").Value = "  -2.08%  "
$ws.Range("E37").Value = "  +1.42%  "
$ws.Range("D38").Value = "'0.01936"
$ws.Range("E38").Value = "  -1.23%  "
$ws.Range("D39").Value = "'0.05141"
$ws.Range("E39").Value = "  -1.30%  "
$ws.Range("D40").Value = "'2.895"
$ws.Range("E40").Value = "  -2.96%  "
$ws.Range("D41").Value = "'6.928"
$ws.Range("E41").Value = "  -1.66%  "
$ws.Range("D42").Value = "'0.4972"
$ws.Range("E42").Value = "  -3.70%  "
$ws.Range("D43").Value = "'0.1601"
$ws.Range("E43").Value = "  -2.32%  "
$ws.Range("D44").Value = "'8.305"
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("D45").Value = "'0.4684"
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("D47").Value = "'10.14"
$ws.Range("E47").Value = "  -1.93%  "
$ws.Range("D48").Value = "'102.24"
$ws.Range("D49").Value = "'1.613"
$ws.Range("E49").Value = "  -2.24%  "
$ws.Range("D50").Value = "'0.06123"
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("D51").Value = "'64.49"
$ws.Range("E51").Value = "  -1.89%  "
